# Fix spelling/content errors across the "Planilha De Requisitos" workbook.
# Order matters: the shared-string table appends newly-introduced text at the
# end in the order cells are edited, so we replicate the exact edit order
# used by the original author to keep the OOXML byte-shape aligned.

$wb = $excel.ActiveWorkbook

$wsRequisitos = $wb.Worksheets.Item("Requisitos")
$wsSprint     = $wb.Worksheets.Item("Sprint Backlog")
$wsBacklog    = $wb.Worksheets.Item("Product Backlog")

# --- Requisitos: non-functional requirement descriptions
# "na residência" (at home) -> "nos ambientes da empresa" / "na empresa" (at the company)
$wsRequisitos.Range("B6").Value  = "Disponibilidade de ar condicionado nos ambientes da empresa"
$wsRequisitos.Range("B8").Value  = "Disponibilidade de umidificador de ar nos ambientes da empresa"
$wsRequisitos.Range("B9").Value  = "Disponibilidade de uma rede de banda larga na empresa"
$wsRequisitos.Range("B11").Value = "Disponibilidade de energia elétrica na empresa"

# --- Sprint Backlog: activity name typos
$wsSprint.Range("E3").Value  = "Cadastrar login e senha no banco de dados"
$wsSprint.Range("E11").Value = "Conexão com o sistema"
$wsSprint.Range("E16").Value = "Conexão com o banco de dados para armazenar os dados recolhidos"
$wsSprint.Range("E20").Value = "Programação de horario de funcionamento da solução"

# --- Product Backlog: description typos
$wsBacklog.Range("D4").Value  = "Automatização do controle da temperatura e umidade através de sensores"
$wsBacklog.Range("D6").Value  = "Recuperação de senha através de email."
$wsBacklog.Range("D10").Value = "Programação opcional de horarios de funcionamento da solução."

# --- Column widths widened to fit the longer corrected text (manual resize,
# no longer auto "best fit")
$wsRequisitos.Columns.Item(2).ColumnWidth = 81.25
$wsSprint.Columns.Item(5).ColumnWidth = 70

# --- Leave the cursor where the author last clicked on each sheet
[void]$wsRequisitos.Range("B14").Select()
[void]$wsSprint.Range("E21").Select()
[void]$wsBacklog.Range("C7").Select()
